$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D and E hold numeric-looking text (prices / percentages) that must
# remain plain text, matching the source data. Force text format before writing
# so Excel does not silently convert them to numeric values.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "302.78"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-1.52%"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-1.63%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.034"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-1.24%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07903"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-2.94%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.845"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-4.40%"
$ws.Range("B7").Value = "KuCoinToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "7.781"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.20%"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9211"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.89%"
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1344"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-2.31%"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1898"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-1.62%"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.09113"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-1.90%"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03466"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-4.07%"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09839"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.19%"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001402"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.28%"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.006134"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "5.92%"
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.715"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "4.54%"
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.107"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-2.00%"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "11.83%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3440"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.08%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "3.11%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.174"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "5.50%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2193"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-8.84%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04412"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-2.33%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001235"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "1.78%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004621"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "5.05%"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "0.14%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01938"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05084"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "2.74%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007612"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-0.47%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.01018"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-8.14%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1343"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-2.79%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002164"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "3.07%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.01017"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-3.81%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006192"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-4.10%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.16%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "63.57"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-1.69%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "39.58%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.16%"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.16%"
